$wb = $excel.ActiveWorkbook

# --- Rename sheets ---
$wsExt = $wb.Worksheets.Item("Test1")
$wsExt.Name = "ExtTest10mm"
$wsFlx = $wb.Worksheets.Item("FlxTest1")
$wsFlx.Name = "FlxTest10mm"

$ws = $wsExt

# --- Update data values on ExtTest10mm sheet ---
# Row 6
$ws.Range("C6").Value = 5.8230000000000004
$ws.Range("D6").Value = 7.3517999999999999
$ws.Range("K6").ClearContents()

# Row 7
$ws.Range("C7").Value = 119
$ws.Range("D7").Value = 114
$ws.Range("K7").ClearContents()

# Row 8
$ws.Range("C8").Formula = "=90-54.6"
$ws.Range("D8").Value = 33.5

# Row 9
$ws.Range("C9").Value = 33.5
$ws.Range("D9").Value = 35.5
$ws.Range("K9").ClearContents()

# Row 10
$ws.Range("C10").Value = 503
$ws.Range("D10").Value = 497

# Row 13
$ws.Range("C13").Value = 40
$ws.Range("D13").Value = 40

# --- Apply yellow fill formatting to row 12 (A12:Q12) ---
$ws.Range("A12:Q12").Interior.Color = 65535

# --- Update selection on the sheet ---
$ws.Range("D6").Select()
